$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 470; this pushes the existing rows 470-515
# down to 471-516 (matching the dimension change A1:R515 -> A1:R516).
$ws.Rows(470).Insert()

# Populate the newly inserted row 470 with the new weekly record.
$ws.Range("A470").Value = 5
$ws.Range("B470").Value = "Macroferia Regional de Talca"
$ws.Range("C470").Value = "Maule"
$ws.Range("D470").Value = 45132
$ws.Range("E470").Value = 7
$ws.Range("F470").Value = 100112006
$ws.Range("G470").Value = "Repollo"
$ws.Range("H470").Value = "Crespo record"
$ws.Range("I470").Value = "Primera"
$ws.Range("J470").Value = 5000
$ws.Range("K470").Value = 600
$ws.Range("L470").Value = 600
$ws.Range("M470").Value = 600
$ws.Range("N470").Value = "$/unidad"
$ws.Range("O470").Value = "Región del Maule"
$ws.Range("P470").Value = 600
$ws.Range("Q470").Value = 1
$ws.Range("R470").Value = "Hortaliza"
